$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.621.02"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +4.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.010.65"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.59%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "508.03"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +8.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.70"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +10.15%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +7.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.64"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +15.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.109"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +13.87%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +8.36%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +5.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.525.52"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.51"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +10.73%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +17.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "56.651.23"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.009.23"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +4.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.81"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +10.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.52"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +10.89%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +12.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "328.37"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +11.32%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +10.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.52"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +7.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.171"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +13.85%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0918"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +15.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.57"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +8.61%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +15.12%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +14.30%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +10.71%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +12.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "156.34"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +13.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.52"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +10.01%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.57%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "24.19"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +5.90%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +10.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.047.98"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +5.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.69"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +8.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.279.31"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +12.26%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +7.42%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +8.88%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +8.14%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +25.36%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +11.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.77"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +8.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.14"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +8.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0872"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +11.72%  "
